$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) - copy the header formatting from G1 (bold,
# bordered, centered) so the new header cell matches the existing ones,
# then set its text.
$g1 = $ws.Range("G1")
$h1 = $ws.Range("H1")
$g1.Copy()
$h1.PasteSpecial(-4122)  # xlPasteFormats
$h1.Value = "Save"

# "Save" flag values for the data rows (2-14)
$saveValues = @(1,0,0,0,0,0,0,1,0,1,0,0,1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
